$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Targeted cell updates (rows 1-177, before the new-row insertion) ---
# Column E ("Reach Rank") holds numeric-looking text in the source file (t="inlineStr");
# a leading apostrophe keeps the COM layer from coercing it to a real number.
$ws.Cells.Item(16,4).Value = 'Restore Reach Function'
$ws.Cells.Item(16,7).Value = 'multiple (HQ pathway)'
$ws.Cells.Item(16,8).Value = 'Cover- Wood'
$ws.Cells.Item(16,9).Value = 'Stability,Coarse Substrate,Flow- Summer Base Flow,Off-Channel- Floodplain,Pool Quantity and Quality,Riparian'
$ws.Cells.Item(27,4).Value = 'Restore Reach Function'
$ws.Cells.Item(27,7).Value = 'multiple (HQ pathway)'
$ws.Cells.Item(31,4).Value = 'Restore Reach Function, Address Limiting Factors'
$ws.Cells.Item(31,7).Value = 'Fry, multiple (HQ pathway)'
$ws.Cells.Item(40,4).Value = 'Restore Reach Function, Address Limiting Factors'
$ws.Cells.Item(40,7).Value = 'Smolt Outmigration, multiple (HQ pathway)'
$ws.Cells.Item(40,8).Value = 'Flow- Summer Base Flow,Riparian,Temperature- Rearing,Predators-Juveniles'
$ws.Cells.Item(40,10).Value = 'Channel Complexity Restoration,Channel Modification,Riparian Restoration and Management,Instream Flow Enhancement,Upland Management,Floodplain Reconnection,Side Channel and Off-Channel Habitat Restoration,Fine Sediment Management,Bank Restoration,Water Quality Improvement,Predator Management'
$ws.Cells.Item(59,5).Value = '''2'
$ws.Cells.Item(68,4).Value = 'Restore Reach Function'
$ws.Cells.Item(68,7).Value = 'multiple (HQ pathway)'
$ws.Cells.Item(72,4).Value = 'Restore Reach Function, Address Limiting Factors'
$ws.Cells.Item(72,7).Value = 'Fry, multiple (HQ pathway)'
$ws.Cells.Item(82,7).Value = 'Holding and Maturation,BT Natal Rearing, multiple (HQ pathway)'
$ws.Cells.Item(86,7).Value = 'Holding and Maturation,Spawning and Incubation,Fry,BT Natal Rearing, multiple (HQ pathway)'
$ws.Cells.Item(120,4).Value = 'Restore Reach Function'
$ws.Cells.Item(120,7).Value = 'multiple (HQ pathway)'
$ws.Cells.Item(124,4).Value = 'Restore Reach Function'
$ws.Cells.Item(124,7).Value = 'multiple (HQ pathway)'
$ws.Cells.Item(124,8).Value = 'Cover- Wood,Riparian,Temperature- Rearing'
$ws.Cells.Item(133,4).Value = 'Restore Reach Function'
$ws.Cells.Item(133,7).Value = 'multiple (HQ pathway)'
$ws.Cells.Item(140,4).Value = 'Restore Reach Function'
$ws.Cells.Item(140,7).Value = 'multiple (HQ pathway)'
$ws.Cells.Item(140,8).Value = 'Cover- Wood,Flow- Summer Base Flow,Off-Channel- Floodplain,Off-Channel- Side-Channels,Riparian'
$ws.Cells.Item(149,7).Value = 'Summer Rearing,Winter Rearing, multiple (HQ pathway)'
$ws.Cells.Item(149,9).Value = 'Coarse Substrate,Pool Quantity and Quality,Cover- Undercut Banks,PRCNT Fines and Embeddedness'
$ws.Cells.Item(157,7).Value = 'Adult Migration,Holding and Maturation,Summer Rearing, multiple (HQ pathway)'
$ws.Cells.Item(161,5).Value = '''1'
$ws.Cells.Item(163,5).Value = '''2'
$ws.Cells.Item(164,5).Value = '''3'
$ws.Cells.Item(165,5).Value = '''3'
$ws.Cells.Item(176,5).Value = '''2'
$ws.Cells.Item(177,5).Value = '''3'

# --- Insert new row 178 ("Methow River Alta Coulee 06") and shift 178-192 down to 179-193 ---
$ws.Rows.Item(178).Insert()

$ws.Cells.Item(178,1).Value = 'Methow River Alta Coulee 06'
$ws.Cells.Item(178,2).Value = 'Methow'
$ws.Cells.Item(178,3).Value = 'Methow River-Alta Coulee'
$ws.Cells.Item(178,4).Value = 'Address Limiting Factors'
$ws.Cells.Item(178,5).Value = '''3'
$ws.Cells.Item(178,6).Value = 'Steelhead'
$ws.Cells.Item(178,7).Value = 'Fry'
$ws.Cells.Item(178,8).Value = 'Cover- Wood'
$ws.Cells.Item(178,9).Value = 'Coarse Substrate,Off-Channel- Floodplain,Off-Channel- Side-Channels'
$ws.Cells.Item(178,10).Value = 'Channel Complexity Restoration,Channel Modification,Fine Sediment Management,Upland Management,Riparian Restoration and Management,Floodplain Reconnection,Side Channel and Off-Channel Habitat Restoration'

Write-Host "Done. UsedRange rows:" $ws.UsedRange.Rows.Count
